$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (day/id headers for the first four subject columns)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON)
$ws.Range("B2").Value = 0.73417456762561273
$ws.Range("C2").Value = 1.9785532574580063
$ws.Range("D2").Value = 0.4871534002765382
$ws.Range("E2").Value = 1.0513439148655763

# Row 3 (STR)
$ws.Range("B3").Value = 0.94891394799480833
$ws.Range("C3").Value = 1.4806712274099791
$ws.Range("D3").Value = 0.92147530980198122
$ws.Range("E3").Value = 0.95963034906768585

# Selection now only spans the edited columns (B:E) instead of the full B:AY range
$ws.Range("B1:E3").Select()
